$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20, shifting rows 20:107 down to 21:108
$ws.Rows.Item(20).Insert()

# Populate the new row 20 with the new data entry
$ws.Cells.Item(20, 1).Value = 3
$ws.Cells.Item(20, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(20, 3).Value = "Coquimbo"
$ws.Cells.Item(20, 4).Value = 44558
$ws.Cells.Item(20, 5).Value = 5
$ws.Cells.Item(20, 6).Value = 100112052
$ws.Cells.Item(20, 7).Value = "Albahaca"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 125
$ws.Cells.Item(20, 11).Value = 4000
$ws.Cells.Item(20, 12).Value = 4500
$ws.Cells.Item(20, 13).Value = 4240
$ws.Cells.Item(20, 14).Value = "`$/docena de matas"
$ws.Cells.Item(20, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(20, 16).Value = 707
$ws.Cells.Item(20, 17).Value = 6
$ws.Cells.Item(20, 18).Value = "Hortaliza"
